# Regenerate the lattice-multiplication exercise table with a new set of
# problems/partial-product digits, cell by cell (table shape is unchanged:
# 5 rows x 3 columns). [char]11 is Word's manual line-break (w:br) char,
# matching Range.Text's representation of line breaks within a cell.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "55 x 84" + [char]11 + "  8    4" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "5|    |"
$t.Cell(1, 2).Range.Text = "33 x 64" + [char]11 + "  6    4" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
$t.Cell(1, 3).Range.Text = "67 x 93" + [char]11 + "  9    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
$t.Cell(2, 1).Range.Text = "66 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(2, 2).Range.Text = "57 x 75" + [char]11 + "  7    5" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "7|    |"
$t.Cell(2, 3).Range.Text = "81 x 91" + [char]11 + "  9    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(3, 1).Range.Text = "22 x 16" + [char]11 + "  1    6" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "2|    |"
$t.Cell(3, 2).Range.Text = "96 x 16" + [char]11 + "  1    6" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "6|    |"
$t.Cell(3, 3).Range.Text = "87 x 28" + [char]11 + "  2    8" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "7|    |"
$t.Cell(4, 1).Range.Text = "16 x 79" + [char]11 + "  7    9" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"
$t.Cell(4, 2).Range.Text = "96 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "6|    |"
$t.Cell(4, 3).Range.Text = "41 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "1|    |"
$t.Cell(5, 1).Range.Text = "40 x 21" + [char]11 + "  2    1" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"
$t.Cell(5, 2).Range.Text = "65 x 34" + [char]11 + "  3    4" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "5|    |"
$t.Cell(5, 3).Range.Text = "67 x 73" + [char]11 + "  7    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
